$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 currently holds the shared string "FETCHTICKET1"; replace its content
# with a single space, matching the author's edit.
$ws.Range("A2").Value = " "

# Rows 8,9,10,12,13,14,16,17 in column A get the same "right aligned,
# no wrap" formatting that column A already uses for rows 2-6 (style
# index 3 in the original workbook). Copy that formatting across instead
# of poking WrapText/HorizontalAlignment individually so we don't leave
# behind any unused intermediate cell-style records.
$ws.Range("A2").Copy()
foreach ($r in @(8, 9, 10, 12, 13, 14, 16, 17)) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

# Update the saved selection to match the workbook as last left by the author.
$ws.Range("A16:A17").Select()
